# "Add Project File (.docx)" -- populate the "Qui ?" (Who?) column (C) of the
# to-do list with the team members assigned to each task, and leave the
# selection on C4 (matching the saved cursor position in the authored file).
#
# New shared strings must be introduced in the order: Alex, Fait, Vicky, Nico
# so they land at shared-string indices 10, 11, 12, 13 respectively (matching
# the authored workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Alex"   # row 8 - "Maquette ?"
$ws.Range("C9").Value = "Alex"   # row 9 - "Use Case"
$ws.Range("C5").Value = "Fait"   # row 5 - "Fonctionnalités"
$ws.Range("C4").Value = "Vicky"  # row 4 - "MCD"
$ws.Range("C3").Value = "Nico"   # row 3 - "Contexte (Entreprise, application)"

# Leave the cursor on C4, matching the workbook's saved selection.
$ws.Range("C4").Select()
